$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-27 Friday" "2026-02-28 Saturday"

Replace-Text "32×59=1888" "43×85=3655"
Replace-Text "91×69=6279" "34×58=1972"
Replace-Text "46×31=1426" "61×13=793"
Replace-Text "77×89=6853" "71×66=4686"
Replace-Text "39×62=2418" "54×60=3240"
Replace-Text "68×87=5916" "72×89=6408"
Replace-Text "68×45=3060" "80×33=2640"
Replace-Text "26×92=2392" "73×63=4599"
Replace-Text "36×78=2808" "94×70=6580"
Replace-Text "64×41=2624" "59×58=3422"
Replace-Text "74×30=2220" "25×50=1250"
Replace-Text "41×81=3321" "98×18=1764"
Replace-Text "86×54=4644" "57×12=684"
Replace-Text "66×93=6138" "12×41=492"
Replace-Text "47×65=3055" "33×70=2310"
Replace-Text "11×68=748" "24×37=888"
Replace-Text "65×75=4875" "30×30=900"
Replace-Text "57×25=1425" "11×60=660"
Replace-Text "31×45=1395" "97×26=2522"
Replace-Text "45×74=3330" "52×41=2132"
Replace-Text "53×64=3392" "32×77=2464"
Replace-Text "62×51=3162" "67×40=2680"
Replace-Text "42×95=3990" "12×11=132"
Replace-Text "23×53=1219" "46×91=4186"
Replace-Text "50×26=1300" "52×41=2132"
